$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename first header (was "Unnamed: 0") to "name"
$ws.Range("B1").Value = "name"

# Insert 3 new blank columns before the current "structure_curated" column (D),
# shifting structure_curated -> G and substance_type_name -> H
$ws.Range("D1:F1").EntireColumn.Insert()

# Set the new header labels
$ws.Range("D1").Value = "meta1"
$ws.Range("E1").Value = "meta2"
$ws.Range("F1").Value = "meta3"

# Fill the new meta columns with dummy placeholder values for every data row
$ws.Range("D2:F21").Value = "dummy"

# Fix up a handful of curated-structure values (column G) that the
# now-working salt remover produces differently
$ws.Range("G9").Value = "[K+]"
$ws.Range("G13").Value = "Nc1ccc2c([O-])c(N=Nc3ccc(-c4ccc(N=Nc5c(S(=O)(=O)O)cc6cc(N)ccc6c5[O-])cc4)cc3)c(S(=O)(=O)O)cc2c1"
$ws.Range("G18").Value = "O=S(=O)([O-])[O-]"
$ws.Range("G19").Value = "[Sb+3]"
